$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with two new columns, P and Q, continuing the
# 0..13 sequence already in B1:O1 (same bold/centered/bordered style).
foreach ($addr in @("P1", "Q1")) {
    $c = $ws.Range($addr)
    $c.Borders.LineStyle = 1
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4160
}
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# For each data row (2-25): flip values in columns I, K, M, O and populate
# the two new columns, P and Q, with 2.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P = 2
    $ws.Cells.Item($r, 17).Value = 2  # Q = 2
}
